$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '61.166.71'
$ws.Cells.Item(2, 5).Value = '  -2.55%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.006.45'
$ws.Cells.Item(3, 5).Value = '  -1.68%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '535.22'
$ws.Cells.Item(5, 5).Value = '  -0.24%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '135.02'
$ws.Cells.Item(6, 5).Value = '  +1.91%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '2.996.85'
$ws.Cells.Item(8, 5).Value = '  -1.72%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.497'
$ws.Cells.Item(9, 5).Value = '  +0.91%  '
$ws.Cells.Item(10, 5).Value = '  -3.21%  '
$ws.Cells.Item(11, 5).Value = '  +0.82%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.448'
$ws.Cells.Item(12, 5).Value = '  -0.49%  '
$ws.Cells.Item(13, 5).Value = '  -0.72%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '34.26'
$ws.Cells.Item(14, 5).Value = '  +0.49%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.494.22'
$ws.Cells.Item(15, 5).Value = '  -1.65%  '
$ws.Cells.Item(16, 5).Value = '  -0.45%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '61.254.76'
$ws.Cells.Item(17, 5).Value = '  -2.45%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '3.002.96'
$ws.Cells.Item(18, 5).Value = '  -1.92%  '
$ws.Cells.Item(19, 5).Value = '  +0.27%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '468.04'
$ws.Cells.Item(20, 5).Value = '  -2.62%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.29'
$ws.Cells.Item(21, 5).Value = '  +0.29%  '
$ws.Cells.Item(22, 5).Value = '  -2.03%  '
$ws.Cells.Item(23, 5).Value = '  -1.48%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '79.75'
$ws.Cells.Item(24, 5).Value = '  +1.09%  '
$ws.Cells.Item(25, 5).Value = '  +0.45%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.15%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.69'
$ws.Cells.Item(27, 5).Value = '  -0.25%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.95'
$ws.Cells.Item(28, 5).Value = '  -1.08%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  +0.04%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.90'
$ws.Cells.Item(30, 5).Value = '  +2.14%  '
$ws.Cells.Item(31, 5).Value = '  +4.17%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '25.57'
$ws.Cells.Item(32, 5).Value = '  -1.17%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.53'
$ws.Cells.Item(33, 5).Value = '  +4.29%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '55.65'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.30'
$ws.Cells.Item(35, 5).Value = '  -2.41%  '
$ws.Cells.Item(36, 5).Value = '  -1.08%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '460.95'
$ws.Cells.Item(37, 5).Value = '  -4.28%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.206.45'
$ws.Cells.Item(38, 5).Value = '  +3.61%  '
$ws.Cells.Item(39, 5).Value = '  -0.19%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.0386'
$ws.Cells.Item(41, 5).Value = '  +2.77%  '
$ws.Cells.Item(42, 5).Value = '  +1.53%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '27.91'
$ws.Cells.Item(43, 5).Value = '  +15.08%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.49'
$ws.Cells.Item(44, 5).Value = '  -4.51%  '
$ws.Cells.Item(45, 5).Value = '  +0.09%  '
$ws.Cells.Item(46, 5).Value = '  -1.36%  '
$ws.Cells.Item(47, 5).Value = '  +0.45%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '120.33'
$ws.Cells.Item(48, 5).Value = '  -0.71%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.108'
$ws.Cells.Item(49, 5).Value = '  +0.69%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0₃0497'
$ws.Cells.Item(50, 5).Value = '  -6.66%  '
$ws.Cells.Item(51, 5).Value = '  +7.96%  '
